$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# ALC row 2
$wsALC.Range("H2").Value = 501
$wsALC.Range("I2").Value = 501
$wsALC.Range("K2").Value = 501
$wsALC.Range("M2").Value = -388

# ALC row 40
$wsALC.Range("H40").Value = 1267.4117
$wsALC.Range("J40").Value = 1233.1111
$wsALC.Range("L40").Value = 1233.1111
$wsALC.Range("N40").Value = -1583.1111

# ALC row 55
$wsALC.Range("H55").Value = 482.29413
$wsALC.Range("I55").Value = 609.9
$wsALC.Range("K55").Value = 609.9
$wsALC.Range("M55").Value = -395.9

# ALC row 62
$wsALC.Range("H62").Value = 0
$wsALC.Range("I62").Value = 0
$wsALC.Range("K62").Value = 0
$wsALC.Range("M62").ClearContents()

# ALC row 65
$wsALC.Range("H65").Value = 0
$wsALC.Range("I65").Value = 0
$wsALC.Range("K65").Value = 0
$wsALC.Range("M65").ClearContents()

# ALC row 92
$wsALC.Range("H92").Value = 3864.7778
$wsALC.Range("I92").Value = 3464.6667
$wsALC.Range("K92").Value = 3464.6667
$wsALC.Range("M92").Value = -2216.6667

# ARM row 2
$wsARM.Range("H2").Value = 806.875
$wsARM.Range("I2").Value = 636.5714
$wsARM.Range("K2").Value = 636.5714
$wsARM.Range("M2").Value = -523.5714

# ARM row 44
$wsARM.Range("H44").Value = 79949
$wsARM.Range("J44").Value = 79949
$wsARM.Range("L44").Value = 79949
$wsARM.Range("N44").Value = -80925

# ARM row 54
$wsARM.Range("H54").Value = 0
$wsARM.Range("J54").Value = 0
$wsARM.Range("L54").Value = 0
$wsARM.Range("N54").ClearContents()

# ARM row 105
$wsARM.Range("H105").Value = 0
$wsARM.Range("J105").Value = 0
$wsARM.Range("L105").Value = 0
$wsARM.Range("N105").ClearContents()

# ARM row 116
$wsARM.Range("H116").Value = 806.875
$wsARM.Range("I116").Value = 636.5714
$wsARM.Range("K116").Value = 636.5714
$wsARM.Range("M116").Value = 1657.4286

# ARM row 122
$wsARM.Range("H122").Value = 9083.615
$wsARM.Range("I122").Value = 9049.125
$wsARM.Range("K122").Value = 27147.375
$wsARM.Range("M122").Value = -24697.375

# BSM row 3
$wsBSM.Range("H3").Value = 806.875
$wsBSM.Range("I3").Value = 636.5714
$wsBSM.Range("K3").Value = 636.5714
$wsBSM.Range("M3").Value = -522.5714

# BSM row 94
$wsBSM.Range("H94").Value = 2454.5454
$wsBSM.Range("I94").Value = 2000
$wsBSM.Range("J94").Value = 2833.3333
$wsBSM.Range("K94").Value = 2000
$wsBSM.Range("L94").Value = 2833.3333
$wsBSM.Range("M94").Value = -1549
$wsBSM.Range("N94").Value = -3735.3333

# CRP row 39
$wsCRP.Range("H39").Value = 18694.334
$wsCRP.Range("I39").Value = 1249.8
$wsCRP.Range("J39").Value = 40500
$wsCRP.Range("K39").Value = 1249.8
$wsCRP.Range("L39").Value = 40500
$wsCRP.Range("M39").Value = -858.8
$wsCRP.Range("N39").Value = -41282

# CRP row 49
$wsCRP.Range("H49").Value = 18694.334
$wsCRP.Range("I49").Value = 1249.8
$wsCRP.Range("J49").Value = 40500
$wsCRP.Range("K49").Value = 1249.8
$wsCRP.Range("L49").Value = 40500
$wsCRP.Range("M49").Value = -1067.8
$wsCRP.Range("N49").Value = -40864

# CRP row 122
$wsCRP.Range("H122").Value = 5793.5
$wsCRP.Range("I122").Value = 399
$wsCRP.Range("J122").Value = 7591.6665
$wsCRP.Range("K122").Value = 1197
$wsCRP.Range("L122").Value = 22774.9995
$wsCRP.Range("M122").Value = 1253
$wsCRP.Range("N122").Value = -27674.9995

# CRP row 124
$wsCRP.Range("H124").Value = 49829.89
$wsCRP.Range("J124").Value = 49829.89
$wsCRP.Range("L124").Value = 49829.89
$wsCRP.Range("N124").Value = -54739.89

# CUL row 42
$wsCUL.Range("H42").Value = 1000
$wsCUL.Range("J42").Value = 1000
$wsCUL.Range("L42").Value = 3000
$wsCUL.Range("N42").Value = -4068

# CUL row 52
$wsCUL.Range("H52").Value = 0
$wsCUL.Range("J52").Value = 0
$wsCUL.Range("L52").Value = 0
$wsCUL.Range("N52").ClearContents()

# CUL row 109
$wsCUL.Range("H109").Value = 244.25
$wsCUL.Range("I109").Value = 264.85715
$wsCUL.Range("J109").Value = 100
$wsCUL.Range("K109").Value = 794.5714499999999
$wsCUL.Range("L109").Value = 300
$wsCUL.Range("M109").Value = 245.4285500000001
$wsCUL.Range("N109").Value = -2380

# GSM row 113
$wsGSM.Range("H113").Value = 3037
$wsGSM.Range("I113").Value = 1999
$wsGSM.Range("J113").Value = 4075
$wsGSM.Range("K113").Value = 1999
$wsGSM.Range("L113").Value = 4075
$wsGSM.Range("M113").Value = 171
$wsGSM.Range("N113").Value = -8415

# GSM row 122
$wsGSM.Range("H122").Value = 5609.8887
$wsGSM.Range("I122").Value = 4399.4
$wsGSM.Range("J122").Value = 7123
$wsGSM.Range("K122").Value = 13198.2
$wsGSM.Range("L122").Value = 21369
$wsGSM.Range("M122").Value = -10748.2
$wsGSM.Range("N122").Value = -26269

# GSM row 126
$wsGSM.Range("H126").Value = 3249.75
$wsGSM.Range("I126").Value = 2666.3333
$wsGSM.Range("J126").Value = 5000
$wsGSM.Range("K126").Value = 7998.999899999999
$wsGSM.Range("L126").Value = 15000
$wsGSM.Range("M126").Value = -5528.999899999999
$wsGSM.Range("N126").Value = -19940

# GSM row 132
$wsGSM.Range("H132").Value = 2948.4167
$wsGSM.Range("I132").Value = 2042.4445
$wsGSM.Range("K132").Value = 6127.333500000001
$wsGSM.Range("M132").Value = -3597.333500000001

# LTW row 40
$wsLTW.Range("H40").Value = 1252996.8
$wsLTW.Range("I40").Value = 3995.6667
$wsLTW.Range("K40").Value = 3995.6667
$wsLTW.Range("M40").Value = -3859.6667

# LTW row 61
$wsLTW.Range("H61").Value = 1697.5
$wsLTW.Range("I61").Value = 1697
$wsLTW.Range("K61").Value = 1697
$wsLTW.Range("M61").Value = -1495

# LTW row 113
$wsLTW.Range("H113").Value = 1697.5
$wsLTW.Range("I113").Value = 1697
$wsLTW.Range("K113").Value = 1697
$wsLTW.Range("M113").Value = 473

# LTW row 132
$wsLTW.Range("H132").Value = 4390.9165
$wsLTW.Range("I132").Value = 7002
$wsLTW.Range("J132").Value = 3520.5557
$wsLTW.Range("K132").Value = 21006
$wsLTW.Range("L132").Value = 10561.6671
$wsLTW.Range("M132").Value = -18476
$wsLTW.Range("N132").Value = -15621.6671

# WVR row 96
$wsWVR.Range("H96").Value = 1150
$wsWVR.Range("I96").Value = 1150
$wsWVR.Range("J96").Value = 0
$wsWVR.Range("K96").Value = 1150
$wsWVR.Range("L96").Value = 0
$wsWVR.Range("M96").Value = 223
$wsWVR.Range("N96").ClearContents()

# WVR row 107
$wsWVR.Range("H107").Value = 567.5714
$wsWVR.Range("I107").Value = 184.9
$wsWVR.Range("K107").Value = 554.7
$wsWVR.Range("M107").Value = 1365.3

# WVR row 126
$wsWVR.Range("H126").Value = 3077.7144
$wsWVR.Range("I126").Value = 3090.6667
$wsWVR.Range("J126").Value = 3000
$wsWVR.Range("K126").Value = 9272.000100000001
$wsWVR.Range("L126").Value = 9000
$wsWVR.Range("M126").Value = -6802.000100000001
$wsWVR.Range("N126").Value = -13940

# WVR row 132
$wsWVR.Range("H132").Value = 2367.4194
$wsWVR.Range("I132").Value = 2396.5667
$wsWVR.Range("J132").Value = 1493
$wsWVR.Range("K132").Value = 7189.7001
$wsWVR.Range("L132").Value = 4479
$wsWVR.Range("M132").Value = -4659.7001
$wsWVR.Range("N132").Value = -9539
